$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2:N6").Value = 85.92117485762657
